$d = $word.ActiveDocument

# Step 1: merge the original 9 runs of the sentence into a single run by
# deleting the whole paragraph's text and retyping it in one shot. This
# triggers the engine's run-consolidation normalization on save.
$full = "If you encounter a merge-conflict error on step 7 in the `"Forking, Cloning and Syncing`" guide while trying to update the contents of master repo in your local directory, executing `"git reset --hard`" before `"git merge upstream/master`" seems to resolve the issue."
$r2 = $d.Range(0, 261)
$r2.Delete()
$ins = $d.Range(0, 0)
$ins.InsertBefore($full)

# Step 2: insert the new "HELLO. " text before everything else.
$ins2 = $d.Range(0, 0)
$ins2.InsertBefore("HELLO. ")

# Step 3: move the _GoBack bookmark so that it sits between the new
# "HELLO. " run and the (merged) sentence run, i.e. right after position 7.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()
$bmRange = $d.Range(7, 7)
$d.Bookmarks.Add("_GoBack", $bmRange)
